# Add a new "Comments" header column (column E) to the four history
# sheets: Withdraw History, Deposit History, Transfer History and
# Absolute History. The "Amount" sheet is left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Withdraw History")
$ws.Range("E1").Value = "Comments"
$ws.Range("E1").Select()

$ws = $wb.Worksheets.Item("Deposit History")
$ws.Range("E1").Value = "Comments"
$ws.Range("E1").Select()

$ws = $wb.Worksheets.Item("Transfer History")
$ws.Range("E1").Value = "Comments"
$ws.Range("E1").Select()

$ws = $wb.Worksheets.Item("Absolute History")
$ws.Range("E1").Value = "Comments"
$ws.Range("E2").Select()
